# ---------------------------------------------------------------------------
# Edit: rename B1 headers on the two existing sheets, and add a new
# "PO Forecast" worksheet (ds / PO_Forecast / yhat_lower / yhat_upper)
# with 67 rows of forecast data, matching the commit's xlsx diff.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1) Rename the "Requested quantity" header on "Weekly Quantity" ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2) Rename the "Requested quantity" header on "Monthly Trend" ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3) Add the new "PO Forecast" worksheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row, styled like the other sheets' header rows (bold, bordered,
# centered horizontally, top-aligned vertically).
$wsForecast.Cells.Item(1,1).Value = "ds"
$wsForecast.Cells.Item(1,2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1,3).Value = "yhat_lower"
$wsForecast.Cells.Item(1,4).Value = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Forecast data rows (dates in column A, integers in B, bounds in C/D).
$data = @(
    @(44934.99999999999,145,-500.7677962670133,788.6562740125271),
    @(44962.99999999999,169,-438.7155412751903,769.9282393035786),
    @(44969.99999999999,175,-484.2055600291505,854.1493955883936),
    @(44997.99999999999,199,-431.7367979564793,850.0654662521536),
    @(45004.99999999999,205,-407.5971696124267,874.4109181504057),
    @(45011.99999999999,211,-412.4089568797188,820.7207156626009),
    @(45025.99999999999,223,-373.6194209032026,870.8315577982767),
    @(45032.99999999999,229,-354.0138622432213,892.0367101067736),
    @(45039.99999999999,235,-340.8248475839889,870.8448617881958),
    @(45046.99999999999,241,-345.887460704431,892.290368620776),
    @(45053.99999999999,247,-354.8043723941048,903.8538766158891),
    @(45060.99999999999,253,-391.0967632147011,881.6912654737064),
    @(45067.99999999999,259,-346.0661794857824,913.8998049653206),
    @(45088.99999999999,277,-360.7956964493622,927.9186672537063),
    @(45109.99999999999,295,-329.5310727953462,952.5149335972839),
    @(45137.99999999999,319,-286.4501008354696,944.8096619100243),
    @(45144.99999999999,325,-318.1586922773513,957.917602354607),
    @(45151.99999999999,331,-268.8303996776658,975.9335303674883),
    @(45158.99999999999,337,-291.826335620028,985.2361604589684),
    @(45165.99999999999,343,-263.3056746737724,1027.658020467412),
    @(45179.99999999999,355,-312.4005421097748,972.5321836867513),
    @(45186.99999999999,361,-277.4278901512515,999.7379639270051),
    @(45193.99999999999,368,-338.3987727135795,985.3825331013586),
    @(45200.99999999999,374,-246.6233170863581,1029.276479279612),
    @(45207.99999999999,380,-177.8729362604974,1024.788239117378),
    @(45214.99999999999,386,-237.0614968614954,1031.069809195132),
    @(45221.99999999999,392,-286.0264976925156,1038.515298261035),
    @(45228.99999999999,398,-246.1944111579259,994.7892279526835),
    @(45242.99999999999,410,-172.7095536586696,1098.914673126329),
    @(45249.99999999999,416,-213.1843678424544,1088.508589792698),
    @(45256.99999999999,422,-221.0841121732408,1005.184315696313),
    @(45263.99999999999,428,-171.7120174849763,1030.218620025457),
    @(45270.99999999999,434,-215.6034774426821,1058.343922451854),
    @(45298.99999999999,458,-153.6822639476281,1041.064242919036),
    @(45305.99999999999,464,-168.4690783882716,1098.329392040604),
    @(45312.99999999999,470,-136.0287066215668,1101.876707177926),
    @(45333.99999999999,488,-127.568626084445,1129.763751343821),
    @(45340.99999999999,494,-156.2068706619311,1107.948620382247),
    @(45347.99999999999,500,-116.5412785301886,1136.531745059375),
    @(45354.99999999999,506,-122.9434305097421,1130.234981480149),
    @(45361.99999999999,512,-125.5959955208131,1111.29033195754),
    @(45396.99999999999,542,-120.2418172160425,1194.453262487953),
    @(45403.99999999999,548,-116.4610076177718,1127.992041398697),
    @(45424.99999999999,566,-50.01748712423286,1148.225917963648),
    @(45431.99999999999,572,-69.19312105417256,1183.459442174422),
    @(45459.99999999999,596,-16.20332025018722,1214.614954098477),
    @(45466.99999999999,602,-17.47630848702078,1214.016353586795),
    @(45473.99999999999,608,-1.530735186970352,1251.758453055968),
    @(45487.99999999999,620,-51.25423689942264,1258.133664955974),
    @(45494.99999999999,626,23.21512282517674,1220.948667084749),
    @(45501.99999999999,632,39.35099303367649,1220.994280154817),
    @(45522.99999999999,650,2.873857220762909,1291.75837292864),
    @(45529.99999999999,656,18.3356064541907,1301.391120077184),
    @(45550.99999999999,674,7.995283143096022,1290.865647810075),
    @(45557.99999999999,680,51.78127857405276,1272.395497945312),
    @(45564.99999999999,686,21.1390865044731,1267.057401110992),
    @(45578.99999999999,699,53.64209056607006,1347.156752663944),
    @(45585.99999999999,705,114.6669531092474,1348.787836924965),
    @(45592.99999999999,711,113.1496505167602,1368.357319380273),
    @(45599.99999999999,717,99.6643315434322,1354.037430220066),
    @(45606.99999999999,723,90.1448858880475,1384.278742391003),
    @(45613.99999999999,729,107.5470547501679,1357.397311382612),
    @(45620.99999999999,735,87.27172577518694,1367.736538207226),
    @(45627.99999999999,741,118.1373359228475,1402.294015876231),
    @(45634.99999999999,747,107.0957025271351,1410.856289991771),
    @(45641.99999999999,753,147.041212533441,1402.413037780157),
    @(45648.99999999999,759,95.12825100562168,1394.602739954754)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# Column A on the data rows uses the same date/time display format as the
# other two sheets' "Order Week"/"Order Month" columns.
$wsForecast.Range("A2:A68").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$wsForecast.Activate()
